# Trabajando en logica de rutas
# Add a new route-log entry to the "rutas_registros" sheet: row 7 gets a
# new date/route pair, and the sheet's used range grows to include the
# following (empty) row 8, mirroring the trailing blank row the sheet
# already carried before this edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rutas_registros")

$ws.Range("A7").Value = "20240720"
$ws.Range("B7").Value = "ruta de test"

# Touch row 8 (without writing real content/formatting) so the sheet's
# used range/dimension extends to include it as a trailing blank row,
# just like the original trailing blank row 7 before this edit.
$ws.Range("A8").WrapText = $false
